$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (shifts Sector..Custom Field 1 right by one)
$ws.Columns("I").Insert()

# New column I: "Instrument" header, "Stock" values for data rows
$ws.Cells(1, 9).Value2 = "Instrument"
$ws.Cells(2, 9).Value2 = "Stock"
$ws.Cells(3, 9).Value2 = "Stock"

# Old Sector column (now J) values change from "Software" to "Tech"
$ws.Cells(2, 10).Value2 = "Tech"
$ws.Cells(3, 10).Value2 = "Tech"

# Update selection to match target (J2 selected)
$ws.Range("J2").Select()
